# Daily attendance processing - reorder "Recorded By" name lists in column G.
# For every row in column G that holds a comma-separated list of recorders,
# reverse the order of the items in that list (single-item cells are left
# untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.Contains(",")) {
        $parts = $val -split ","
        $n = $parts.Length

        $reversedParts = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i].Trim()
        }

        $cell.Value2 = [string]::Join(", ", $reversedParts)
    }
}
